$wb = $excel.ActiveWorkbook

# ALC!row21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 27049.133
$ws.Range("I21").Value = 32504.75
$ws.Range("J21").Value = 20814.143
$ws.Range("K21").Value = 32504.75
$ws.Range("L21").Value = 20814.143
$ws.Range("M21").Value = -32036.75
$ws.Range("N21").Value = -21750.143

# ALC!row23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 27049.133
$ws.Range("I23").Value = 32504.75
$ws.Range("J23").Value = 20814.143
$ws.Range("K23").Value = 32504.75
$ws.Range("L23").Value = 20814.143
$ws.Range("M23").Value = -32270.75
$ws.Range("N23").Value = -21282.143

# ALC!row62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4663.727
$ws.Range("I62").Value = 6549.1665
$ws.Range("K62").Value = 6549.1665
$ws.Range("M62").Value = -5925.1665

# ALC!row65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4663.727
$ws.Range("I65").Value = 6549.1665
$ws.Range("K65").Value = 32745.8325
$ws.Range("M65").Value = -29625.8325

# ALC!row87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 37175
$ws.Range("J87").Value = 37175
$ws.Range("L87").Value = 37175
$ws.Range("N87").Value = -39671

# ALC!row90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 37175
$ws.Range("J90").Value = 37175
$ws.Range("L90").Value = 111525
$ws.Range("N90").Value = -124005

# ALC!row130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 48502
$ws.Range("J130").Value = 48502
$ws.Range("L130").Value = 48502
$ws.Range("N130").Value = -58542

# ARM!row22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 2800
$ws.Range("J22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("N22").Value = -8598

# ARM!row39
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# ARM!row80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 53317
$ws.Range("J80").Value = 53317
$ws.Range("L80").Value = 53317
$ws.Range("N80").Value = -55313

# ARM!row83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 53317
$ws.Range("J83").Value = 53317
$ws.Range("L83").Value = 159951
$ws.Range("N83").Value = -169935

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 18204.5
$ws.Range("I102").Value = 1855
$ws.Range("J102").Value = 20929.416
$ws.Range("K102").Value = 1855
$ws.Range("L102").Value = 20929.416
$ws.Range("M102").Value = -233
$ws.Range("N102").Value = -24173.416

# ARM!row114
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 36364
$ws.Range("J114").Value = 36364
$ws.Range("L114").Value = 36364
$ws.Range("N114").Value = -45042

# ARM!row119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 52683.668
$ws.Range("J119").Value = 52683.668
$ws.Range("L119").Value = 52683.668
$ws.Range("N119").Value = -62359.668

# ARM!row121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 32495.75
$ws.Range("J121").Value = 32495.75
$ws.Range("L121").Value = 32495.75
$ws.Range("N121").Value = -35989.75

# ARM!row131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 48676.75
$ws.Range("J131").Value = 48676.75
$ws.Range("L131").Value = 48676.75
$ws.Range("N131").Value = -58756.75

# BSM!row108
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 42729.6
$ws.Range("J108").Value = 42729.6
$ws.Range("L108").Value = 42729.6
$ws.Range("N108").Value = -50409.6

# BSM!row112
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 45975
$ws.Range("J112").Value = 45975
$ws.Range("L112").Value = 45975
$ws.Range("N112").Value = -48929

# BSM!row130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 39694.668
$ws.Range("J130").Value = 39694.668
$ws.Range("L130").Value = 39694.668
$ws.Range("N130").Value = -49734.668

# BSM!row132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 54240
$ws.Range("J132").Value = 54240
$ws.Range("L132").Value = 54240
$ws.Range("N132").Value = -64360

# CRP!row64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30545.166
$ws.Range("J64").Value = 30545.166
$ws.Range("L64").Value = 30545.166
$ws.Range("N64").Value = -31041.166

# CRP!row67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 30545.166
$ws.Range("J67").Value = 30545.166
$ws.Range("L67").Value = 30545.166
$ws.Range("N67").Value = -32261.166

# CRP!row111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 47600.332
$ws.Range("J111").Value = 47600.332
$ws.Range("L111").Value = 47600.332
$ws.Range("N111").Value = -55780.332

# CRP!row118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 55988.77
$ws.Range("I132").Value = 1810.5625
$ws.Range("J132").Value = 142673.9
$ws.Range("K132").Value = 5431.6875
$ws.Range("L132").Value = 428021.7
$ws.Range("M132").Value = -2901.6875
$ws.Range("N132").Value = -433081.7

# CRP!row138
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 40200
$ws.Range("J138").Value = 40200
$ws.Range("L138").Value = 40200
$ws.Range("N138").Value = -50480

# CUL!row48
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 334833.34
$ws.Range("J48").Value = 334833.34
$ws.Range("L48").Value = 1004500.02
$ws.Range("N48").Value = -1005000.02

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3388.4866
$ws.Range("I113").Value = 4828.7085
$ws.Range("J113").Value = 729.61536
$ws.Range("K113").Value = 14486.1255
$ws.Range("L113").Value = 2188.84608
$ws.Range("M113").Value = -12316.1255
$ws.Range("N113").Value = -6528.84608

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3250.9363
$ws.Range("J131").Value = 1209.1464
$ws.Range("L131").Value = 3627.4392
$ws.Range("N131").Value = -13707.4392

# GSM!row20
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 60003320
$ws.Range("J20").Value = 4150
$ws.Range("L20").Value = 4150
$ws.Range("N20").Value = -4640

# GSM!row21
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 7840
$ws.Range("J21").Value = 7840
$ws.Range("L21").Value = 7840
$ws.Range("N21").Value = -8186

# GSM!row30
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 7840
$ws.Range("J30").Value = 7840
$ws.Range("L30").Value = 7840
$ws.Range("N30").Value = -8050

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 189222.52
$ws.Range("I80").Value = 337400.6
$ws.Range("J80").Value = 3999.9167
$ws.Range("K80").Value = 337400.6
$ws.Range("L80").Value = 3999.9167
$ws.Range("M80").Value = -336402.6
$ws.Range("N80").Value = -5995.9167

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 189222.52
$ws.Range("I83").Value = 337400.6
$ws.Range("J83").Value = 3999.9167
$ws.Range("K83").Value = 1687003
$ws.Range("L83").Value = 19999.5835
$ws.Range("M83").Value = -1682011
$ws.Range("N83").Value = -29983.5835

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1143.5555
$ws.Range("I122").Value = 1184.5714
$ws.Range("K122").Value = 3553.7142
$ws.Range("M122").Value = -1103.7142

# GSM!row130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 44859.89
$ws.Range("J130").Value = 44859.89
$ws.Range("L130").Value = 44859.89
$ws.Range("N130").Value = -54899.89

# GSM!row138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 37973.332
$ws.Range("J138").Value = 37973.332
$ws.Range("L138").Value = 37973.332
$ws.Range("N138").Value = -48253.332

# LTW!row21
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 39800
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1289.1818
$ws.Range("I22").Value = 2000.5
$ws.Range("J22").Value = 1131.1111
$ws.Range("K22").Value = 2000.5
$ws.Range("L22").Value = 1131.1111
$ws.Range("M22").Value = -1705.5
$ws.Range("N22").Value = -1721.1111

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1289.1818
$ws.Range("I27").Value = 2000.5
$ws.Range("J27").Value = 1131.1111
$ws.Range("K27").Value = 2000.5
$ws.Range("L27").Value = 1131.1111
$ws.Range("M27").Value = -1893.5
$ws.Range("N27").Value = -1345.1111

# LTW!row76
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 14799.5
$ws.Range("J76").Value = 19968.5
$ws.Range("L76").Value = 19968.5
$ws.Range("N76").Value = -20644.5

# LTW!row79
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 14799.5
$ws.Range("J79").Value = 19968.5
$ws.Range("L79").Value = 19968.5
$ws.Range("N79").Value = -22308.5

# LTW!row111
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 40713.6
$ws.Range("J111").Value = 40713.6
$ws.Range("L111").Value = 40713.6
$ws.Range("N111").Value = -48893.6

# LTW!row128
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 43714.5
$ws.Range("J128").Value = 43714.5
$ws.Range("L128").Value = 43714.5
$ws.Range("N128").Value = -53674.5

# WVR!row18
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 9766.667
$ws.Range("I18").Value = 8600
$ws.Range("K18").Value = 8600
$ws.Range("M18").Value = -8427

# WVR!row108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 28430.666
$ws.Range("J108").Value = 28430.666
$ws.Range("L108").Value = 28430.666
$ws.Range("N108").Value = -36110.666

Write-Output "applied all Masamune_Profits updates"
